$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 383. Excel shifts the
# existing rows 383:405 down to 385:407 and copies the row-383 formatting
# (incl. the date style on column D) onto the freshly inserted rows.
$ws.Rows("383:384").Insert()

# New row 383: "Pintón" quality entry for the week of 2021-11-16 (serial 44516)
$ws.Range("A383").Value = 7
$ws.Range("B383").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C383").Value = "Ñuble"
$ws.Range("D383").Value = 44516
$ws.Range("E383").Value = 16
$ws.Range("F383").Value = "Fruta"
$ws.Range("G383").Value = 100108
$ws.Range("H383").Value = "Tropicales y subtropicales"
$ws.Range("I383").Value = 100108006
$ws.Range("J383").Value = "Plátano"
$ws.Range("K383").Value = "Sin especificar"
$ws.Range("L383").Value = "Pintón"
$ws.Range("M383").Value = 180
$ws.Range("N383").Value = 16000
$ws.Range("O383").Value = 16000
$ws.Range("P383").Value = 16000
$ws.Range("Q383").Value = "`$/caja 20 kilos"
$ws.Range("R383").Value = "Ecuador"
$ws.Range("S383").Value = 800
$ws.Range("T383").Value = 20

# New row 384: "Primera Pintón" quality entry for the same week
$ws.Range("A384").Value = 7
$ws.Range("B384").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C384").Value = "Ñuble"
$ws.Range("D384").Value = 44516
$ws.Range("E384").Value = 16
$ws.Range("F384").Value = "Fruta"
$ws.Range("G384").Value = 100108
$ws.Range("H384").Value = "Tropicales y subtropicales"
$ws.Range("I384").Value = 100108006
$ws.Range("J384").Value = "Plátano"
$ws.Range("K384").Value = "Sin especificar"
$ws.Range("L384").Value = "Primera Pintón"
$ws.Range("M384").Value = 240
$ws.Range("N384").Value = 17000
$ws.Range("O384").Value = 18000
$ws.Range("P384").Value = 17500
$ws.Range("Q384").Value = "`$/caja 20 kilos"
$ws.Range("R384").Value = "Ecuador"
$ws.Range("S384").Value = 875
$ws.Range("T384").Value = 20
